# [Excel] Update range edge sample (#549)
# Remove redundant range edge ref docs samples:
#   - "getExtendedRange" / "getExtendedRangeRight" sample row
#   - "getRangeEdge" / "getRangeEdgeLeft" sample row
# These were duplicate/redundant sample rows for the same method, each
# documenting a second (less canonical) edge direction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the higher-numbered row first so the lower row's index doesn't shift
# before we delete it.
$ws.Rows("147:147").Delete()   # getRangeEdge / getRangeEdgeLeft
$ws.Rows("142:142").Delete()   # getExtendedRange / getExtendedRangeRight

# Restore the active selection to where editing left off.
$ws.Range("A146:XFD146").Select()
